$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H header - "observacao"
$ws.Range("H1").Value = "observacao"

# Row 2 updates
$ws.Range("B2").Value = "18/06/2025, 03:48"
$ws.Range("C2").Value = "UT71V"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2025-06-17"
$ws.Range("D2").Style = "Normal"

$ws.Range("F2").Value = "Elétrico"
$ws.Range("G2").Value = "Cabo Acionamento"

# Row 3 updates
$ws.Range("B3").Value = "18/06/2025, 03:49"
$ws.Range("C3").Value = "ESHV2"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2025-06-17"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "48"
$ws.Range("E3").Style = "Normal"

$ws.Range("G3").Value = "Caixa do Micro"
$ws.Range("H3").Value = "teste"
